$d = $word.ActiveDocument

# Locate the paragraph that starts with "Não existem estatísticas..." and the
# paragraph right after it ("Ricardo Augusto Dias...") - together they are being
# replaced / expanded into several new paragraphs of content.
$startPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Não existem estatísticas oficiais*") {
        $startPara = $p
        break
    }
}
if ($startPara -eq $null) {
    throw "Could not locate the Não existem estatísticas... paragraph"
}
$endPara = $startPara.Next()
if ($endPara -eq $null -or -not ($endPara.Range.Text -like "Ricardo Augusto Dias*")) {
    throw "Unexpected document structure: paragraph after target is not Ricardo Augusto Dias..."
}

$targetRange = $d.Range($startPara.Range.Start, $endPara.Range.End)

# NOTE: InsertXML on this host mis-behaves (wipes the whole document body)
# when a single call inserts 5+ paragraphs worth of open-xml, so the new
# content (7 paragraphs) is applied in smaller batches: the first batch
# replaces the old target range, subsequent batches are appended right
# after the content that was just inserted.
$xmlChunk0 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Qual seria então a solução para o problema de maus-tratos de animais?</w:t></w:r><w:r><w:t xml:space="preserve"> Trabalhar nas causas, por meio da castração e a informação da posse responsável. A responsabilidade das pessoas de tratar bem o animal que está adquirindo, e quando for adquirir, seja por compra ou por adoção, ter a certeza de que está tomando a decisão correta, e que há viabilidade para receber o animal na residência. Adquirir um animal não é apenas leva-lo para casa, brincar e dar comida, antes de tomar a decisão de ter um é necessário pensar, analisar, se informar a respeito de raças, cuidados, gastos, espaço, tempo, tamanho, etc. Ser prudente na decisão de ter ou não um animal é a melhor indicação para não colocar em risco a vida dos animais.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Todos que se interessam em obter um animal deveriam ter acesso aos horrores que acontecem a animais que se encontram nas ruas, e acabam muitos deles com um fim trágico. Além de gerarem zoonoses, esses animais têm um destino que </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>ninguém, em sã consciência, desejaria à mais ínfima criatura.</w:t></w:r><w:r><w:t xml:space="preserve"> Eles passam por fome, sede, frio, calor, além dos maus-tratos pelas ruas.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Enquanto há uma mortalidade de 16 mil animais/ano por eutanásia no CCZ do município de São Paulo (dados de 2003) (Pet </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Food</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Health </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>and</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Care</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, n. 4), nas outras cidades do Brasil o número de animais que morrem é igual ou maior. Não é um número assustador? O CCZ de São Paulo tem as instalações comparadas as de países de primeiro mundo, o que deveria implicar em uma diminuição do número de mortes por eutanásia e de animais maltratados, mas nada disso adianta se a população não contribuir.</w:t></w:r><w:r><w:t xml:space="preserve"> Só assim poderemos erradicar o número de animais doentes e sacrificados.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Não existem estatísticas oficiais sobre o número de animais desamparados nas ruas, pelo fato de medir essa quantidade ser uma tarefa extremamente difícil. De acordo levantamento realizado pela VEJA SÃO PAULO, em 10 das principais instituições da </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>capital paulista, cerca de 500 animais são resgatados das ruas por mês, totalizando 6000 por ano. Segundo os profissionais dessas ONGs, grande parte deles já teve um lar. Esse número trata-se apenas de uma amostragem, de acordo com os especialistas o problema que vivemos hoje com relação ao abandono de animais é muito maior.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$targetRange.InsertXML($xmlChunk0)

$xmlChunk1 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Ricardo Augusto Dias, professor da Faculdade de Medicina Veterinária e Zootecnia da Universidade de São Paulo, afirma que os animais de rua</w:t></w:r><w:r><w:t xml:space="preserve"> costumam se concentrar em áreas de limpeza escassa e com abrigo, como terrenos baldios e construções. Além disso, alguns têm endereço fixo, mas contam com acesso à rua, outros estão perdidos e há os chamados “cães comunitários”, cuidados por diversas pessoas.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Os casos de animais que já tiveram um dono e um lar, e hoje viraram “órfãos”, são de cortar o coração. Por mais que a ideia de considerar o animal doméstico como um membro da família esteja se expandindo, muitas pessoas ainda insistem em trata-los como mercadoria, um objeto que pode ser descartado. </w:t></w:r><w:r><w:t>“Já ouvi os motivos mais absurdos de tutores para desistir das mascotes, do naipe de ‘fiquei grávida’ ou ‘comecei a nam</w:t></w:r><w:r><w:t xml:space="preserve">orar e minha parceira tem medo’ </w:t></w:r><w:r><w:t>”, diz a ativista Luisa Mell, cujo instituto recebe cerca de 500 pedidos de resgate diariamente.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Todo fim de ano, o aumento do abandono de animais é notável. Com as festas, muitos optam por viajar e não sabem o que fazer com o animal. Hoje em dia, existem </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">hotéis próprios para receber animais domésticos em casos como esse, porém o custo é alto e muitas pessoas preferem abandonar o animal. </w:t></w:r><w:r><w:t>“Nunca me esqueci de quando fui procurada por uma mulher que ia se mudar de casa e queria deixar comigo seu cachorro de 10 anos. Como pode jog</w:t></w:r><w:r><w:t>ar fora um companheiro de uma dé</w:t></w:r><w:r><w:t>cada? ”, espanta-se Luisa</w:t></w:r><w:r><w:t xml:space="preserve"> Mell</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
# Guard: walk back from doc end in case the last paragraph is the sectPr-only one
$insertPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$insertPoint.InsertXML($xmlChunk1)

